# Update the three-digit x one-digit multiplication problems/answers
# in the table to the new set of values, per the commit's regenerated
# output.

$d = $word.ActiveDocument

$pairs = @(
    @("497×4=1988", "983×2=1966"),
    @("763×8=6104", "257×5=1285"),
    @("177×2=354",  "584×7=4088"),
    @("211×9=1899", "742×6=4452"),
    @("339×2=678",  "280×5=1400"),
    @("565×8=4520", "240×3=720"),
    @("666×9=5994", "786×6=4716"),
    @("747×9=6723", "896×2=1792"),
    @("821×4=3284", "400×3=1200"),
    @("102×4=408",  "820×7=5740"),
    @("751×7=5257", "962×3=2886"),
    @("496×3=1488", "762×5=3810"),
    @("678×6=4068", "856×9=7704"),
    @("899×2=1798", "636×7=4452"),
    @("342×5=1710", "188×9=1692"),
    @("145×7=1015", "211×3=633"),
    @("237×5=1185", "420×9=3780"),
    @("110×9=990",  "364×2=728"),
    @("269×3=807",  "667×4=2668"),
    @("409×6=2454", "383×7=2681"),
    @("450×2=900",  "937×6=5622"),
    @("901×3=2703", "125×4=500"),
    @("623×7=4361", "421×7=2947"),
    @("962×2=1924", "863×4=3452"),
    @("151×8=1208", "313×3=939")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
